# Add a new "usuario" row (Juan Pablo) at the bottom of the sheet, mirroring
# the layout/styling of the existing rows, and wire up a mailto hyperlink on
# the new email cell (same pattern used for the "Parra" row above it).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRow = 10

# Write the new row's values in column order (A -> E) so the shared-string
# table grows in the same order the source workbook uses.
$ws.Cells.Item($newRow, 1).Value = "Juan Pablo"
$ws.Cells.Item($newRow, 2).Value = "juan.villagra2201@alumnos.ubiobio.cl"
$ws.Cells.Item($newRow, 3).Value = "20680787-3"
$ws.Cells.Item($newRow, 4).Value = "ALUMNO"
$ws.Cells.Item($newRow, 5).Value = "Ingeniería Civil en Informática"

# Reuse the formatting of the row above (A9/C9/D9/E9 share style index 4)
# for the new row's Name/RUT/Rol cells.
$ws.Cells.Item($newRow - 1, 1).Copy()
$ws.Cells.Item($newRow, 1).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 3).PasteSpecial(-4122)
$ws.Cells.Item($newRow, 4).PasteSpecial(-4122)

# Turn the email cell into a mailto hyperlink (this also applies the
# workbook's hyperlink cell style automatically).
$ws.Hyperlinks.Add($ws.Range("B" + $newRow), "mailto:juan.villagra2201@alumnos.ubiobio.cl")

# Leave the selection on the last cell of the new row, like Excel would
# after tabbing through the newly entered data.
$ws.Application.CutCopyMode = $false
$ws.Range("E10").Select()
